# Rebase the ipc_inquilinos (C) and ipc_oficial (E) series so that the
# first observation (row 2, 2019-01) equals 100, then propagate the same
# scaling factor down the whole series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 80

$baseC = $ws.Range("C$firstRow").Value2
$baseE = $ws.Range("E$firstRow").Value2

$ratioC = 100.0 / $baseC
$ratioE = 100.0 / $baseE

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cAddr = "C" + $r
    $eAddr = "E" + $r

    $oldC = $ws.Range($cAddr).Value2
    $oldE = $ws.Range($eAddr).Value2

    $ws.Range($cAddr).Value = $oldC * $ratioC
    $ws.Range($eAddr).Value = $oldE * $ratioE
}
